$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 ("Marking"): Right 6 -> 9, Wrong 3 -> 2
$ws.Range("B11").Value = 9
$ws.Range("C11").Value = 2

# Row 12 ("Total"): Right 156 -> 234, fraction label "156/168" -> "234/252"
$ws.Range("B12").Value = 234
$ws.Range("E12").Value = "234/252"
